$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 54; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = [string]$cell.Value2
    if ($val.EndsWith("16")) {
        $cell.Value2 = $val.Substring(0, $val.Length - 2)
    }
}
